# Natmi following Dr Hou advice
# Recomputed LR-pair statistics (Lama3-Sdc2) across all 4x4 sending/target
# cluster combinations (ECs, FAPs, M2, sCs), expanding from 12 to 16 rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Lama3"
$ws.Cells.Item(2, 3).Value = "Sdc2"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 4.736398666666666
$ws.Cells.Item(2, 8).Value = 14.209196
$ws.Cells.Item(2, 9).Value = 0.7599472210591204
$ws.Cells.Item(2, 10).Value = 0.7599472210591204
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 3.057109
$ws.Cells.Item(2, 14).Value = 9.171327
$ws.Cells.Item(2, 15).Value = 0.02694952608666365
$ws.Cells.Item(2, 16).Value = 0.02694952608666365
$ws.Cells.Item(2, 17).Value = 14.47968699145467
$ws.Cells.Item(2, 18).Value = 130.317182923092
$ws.Cells.Item(2, 19).Value = 0.02048021745842031
$ws.Cells.Item(2, 20).Value = 0.02048021745842031

# Row 3: ECs -> FAPs
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Lama3"
$ws.Cells.Item(3, 3).Value = "Sdc2"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 4.736398666666666
$ws.Cells.Item(3, 8).Value = 14.209196
$ws.Cells.Item(3, 9).Value = 0.7599472210591204
$ws.Cells.Item(3, 10).Value = 0.7599472210591204
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 89.02756
$ws.Cells.Item(3, 14).Value = 267.08268
$ws.Cells.Item(3, 15).Value = 0.7848102735793893
$ws.Cells.Item(3, 16).Value = 0.7848102735793893
$ws.Cells.Item(3, 17).Value = 421.6700164805866
$ws.Cells.Item(3, 18).Value = 3795.030148325279
$ws.Cells.Item(3, 19).Value = 0.5964143864653049
$ws.Cells.Item(3, 20).Value = 0.5964143864653049

# Row 4: ECs -> M2
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Lama3"
$ws.Cells.Item(4, 3).Value = "Sdc2"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 4.736398666666666
$ws.Cells.Item(4, 8).Value = 14.209196
$ws.Cells.Item(4, 9).Value = 0.7599472210591204
$ws.Cells.Item(4, 10).Value = 0.7599472210591204
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.184005
$ws.Cells.Item(4, 14).Value = 0.5520149999999999
$ws.Cells.Item(4, 15).Value = 0.0016220709001794
$ws.Cells.Item(4, 16).Value = 0.0016220709001794
$ws.Cells.Item(4, 17).Value = 0.8715210366599998
$ws.Cells.Item(4, 18).Value = 7.843689329939998
$ws.Cells.Item(4, 19).Value = 0.001232688272952201
$ws.Cells.Item(4, 20).Value = 0.001232688272952201

# Row 5: ECs -> sCs
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Lama3"
$ws.Cells.Item(5, 3).Value = "Sdc2"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 4.736398666666666
$ws.Cells.Item(5, 8).Value = 14.209196
$ws.Cells.Item(5, 9).Value = 0.7599472210591204
$ws.Cells.Item(5, 10).Value = 0.7599472210591204
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 21.16964733333333
$ws.Cells.Item(5, 14).Value = 63.508942
$ws.Cells.Item(5, 15).Value = 0.1866181294337677
$ws.Cells.Item(5, 16).Value = 0.1866181294337677
$ws.Cells.Item(5, 17).Value = 100.2678894034035
$ws.Cells.Item(5, 18).Value = 902.4110046306319
$ws.Cells.Item(5, 19).Value = 0.141819928862443
$ws.Cells.Item(5, 20).Value = 0.141819928862443

# Row 6: FAPs -> ECs
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Lama3"
$ws.Cells.Item(6, 3).Value = "Sdc2"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0.2434186666666667
$ws.Cells.Item(6, 8).Value = 0.730256
$ws.Cells.Item(6, 9).Value = 0.03905611674733385
$ws.Cells.Item(6, 10).Value = 0.03905611674733384
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 3.057109
$ws.Cells.Item(6, 14).Value = 9.171327
$ws.Cells.Item(6, 15).Value = 0.02694952608666365
$ws.Cells.Item(6, 16).Value = 0.02694952608666365
$ws.Cells.Item(6, 17).Value = 0.7441573966346667
$ws.Cells.Item(6, 18).Value = 6.697416569712
$ws.Cells.Item(6, 19).Value = 0.001052543837126055
$ws.Cells.Item(6, 20).Value = 0.001052543837126054

# Row 7: FAPs -> FAPs
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Lama3"
$ws.Cells.Item(7, 3).Value = "Sdc2"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.2434186666666667
$ws.Cells.Item(7, 8).Value = 0.730256
$ws.Cells.Item(7, 9).Value = 0.03905611674733385
$ws.Cells.Item(7, 10).Value = 0.03905611674733384
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 89.02756
$ws.Cells.Item(7, 14).Value = 267.08268
$ws.Cells.Item(7, 15).Value = 0.7848102735793893
$ws.Cells.Item(7, 16).Value = 0.7848102735793893
$ws.Cells.Item(7, 17).Value = 21.67096995178667
$ws.Cells.Item(7, 18).Value = 195.03872956608
$ws.Cells.Item(7, 19).Value = 0.03065164166942365
$ws.Cells.Item(7, 20).Value = 0.03065164166942364

# Row 8: FAPs -> M2
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Lama3"
$ws.Cells.Item(8, 3).Value = "Sdc2"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.2434186666666667
$ws.Cells.Item(8, 8).Value = 0.730256
$ws.Cells.Item(8, 9).Value = 0.03905611674733385
$ws.Cells.Item(8, 10).Value = 0.03905611674733384
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.184005
$ws.Cells.Item(8, 14).Value = 0.5520149999999999
$ws.Cells.Item(8, 15).Value = 0.0016220709001794
$ws.Cells.Item(8, 16).Value = 0.0016220709001794
$ws.Cells.Item(8, 17).Value = 0.04479025176
$ws.Cells.Item(8, 18).Value = 0.4031122658399999
$ws.Cells.Item(8, 19).Value = 0.00006335179044985955
$ws.Cells.Item(8, 20).Value = 0.00006335179044985952

# Row 9: FAPs -> sCs
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Lama3"
$ws.Cells.Item(9, 3).Value = "Sdc2"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.2434186666666667
$ws.Cells.Item(9, 8).Value = 0.730256
$ws.Cells.Item(9, 9).Value = 0.03905611674733385
$ws.Cells.Item(9, 10).Value = 0.03905611674733384
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 21.16964733333333
$ws.Cells.Item(9, 14).Value = 63.508942
$ws.Cells.Item(9, 15).Value = 0.1866181294337677
$ws.Cells.Item(9, 16).Value = 0.1866181294337677
$ws.Cells.Item(9, 17).Value = 5.153087327683556
$ws.Cells.Item(9, 18).Value = 46.377785949152
$ws.Cells.Item(9, 19).Value = 0.007288579450334292
$ws.Cells.Item(9, 20).Value = 0.00728857945033429

# Row 10: M2 -> ECs
$ws.Cells.Item(10, 1).Value = "M2"
$ws.Cells.Item(10, 2).Value = "Lama3"
$ws.Cells.Item(10, 3).Value = "Sdc2"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.07790599999999999
$ws.Cells.Item(10, 8).Value = 0.233718
$ws.Cells.Item(10, 9).Value = 0.01249988701763953
$ws.Cells.Item(10, 10).Value = 0.01249988701763953
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 3.057109
$ws.Cells.Item(10, 14).Value = 9.171327
$ws.Cells.Item(10, 15).Value = 0.02694952608666365
$ws.Cells.Item(10, 16).Value = 0.02694952608666365
$ws.Cells.Item(10, 17).Value = 0.238167133754
$ws.Cells.Item(10, 18).Value = 2.143504203786
$ws.Cells.Item(10, 19).Value = 0.0003368660312622247
$ws.Cells.Item(10, 20).Value = 0.0003368660312622247

# Row 11: M2 -> FAPs
$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 2).Value = "Lama3"
$ws.Cells.Item(11, 3).Value = "Sdc2"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 0.07790599999999999
$ws.Cells.Item(11, 8).Value = 0.233718
$ws.Cells.Item(11, 9).Value = 0.01249988701763953
$ws.Cells.Item(11, 10).Value = 0.01249988701763953
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 89.02756
$ws.Cells.Item(11, 14).Value = 267.08268
$ws.Cells.Item(11, 15).Value = 0.7848102735793893
$ws.Cells.Item(11, 16).Value = 0.7848102735793893
$ws.Cells.Item(11, 17).Value = 6.935781089359999
$ws.Cells.Item(11, 18).Value = 62.42202980423999
$ws.Cells.Item(11, 19).Value = 0.009810039750025134
$ws.Cells.Item(11, 20).Value = 0.009810039750025134

# Row 12: M2 -> M2
$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 2).Value = "Lama3"
$ws.Cells.Item(12, 3).Value = "Sdc2"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 0.07790599999999999
$ws.Cells.Item(12, 8).Value = 0.233718
$ws.Cells.Item(12, 9).Value = 0.01249988701763953
$ws.Cells.Item(12, 10).Value = 0.01249988701763953
$ws.Cells.Item(12, 11).Value = 2
$ws.Cells.Item(12, 12).Value = 0.6666666666666666
$ws.Cells.Item(12, 13).Value = 0.184005
$ws.Cells.Item(12, 14).Value = 0.5520149999999999
$ws.Cells.Item(12, 15).Value = 0.0016220709001794
$ws.Cells.Item(12, 16).Value = 0.0016220709001794
$ws.Cells.Item(12, 17).Value = 0.01433509353
$ws.Cells.Item(12, 18).Value = 0.12901584177
$ws.Cells.Item(12, 19).Value = 0.00002027570298684334
$ws.Cells.Item(12, 20).Value = 0.00002027570298684334

# Row 13: M2 -> sCs
$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 2).Value = "Lama3"
$ws.Cells.Item(13, 3).Value = "Sdc2"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 0.07790599999999999
$ws.Cells.Item(13, 8).Value = 0.233718
$ws.Cells.Item(13, 9).Value = 0.01249988701763953
$ws.Cells.Item(13, 10).Value = 0.01249988701763953
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 21.16964733333333
$ws.Cells.Item(13, 14).Value = 63.508942
$ws.Cells.Item(13, 15).Value = 0.1866181294337677
$ws.Cells.Item(13, 16).Value = 0.1866181294337677
$ws.Cells.Item(13, 17).Value = 1.649242545150666
$ws.Cells.Item(13, 18).Value = 14.843182906356
$ws.Cells.Item(13, 19).Value = 0.002332705533365326
$ws.Cells.Item(13, 20).Value = 0.002332705533365326

# Row 14: sCs -> ECs
$ws.Cells.Item(14, 1).Value = "sCs"
$ws.Cells.Item(14, 2).Value = "Lama3"
$ws.Cells.Item(14, 3).Value = "Sdc2"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 1.174813
$ws.Cells.Item(14, 8).Value = 3.524439
$ws.Cells.Item(14, 9).Value = 0.1884967751759062
$ws.Cells.Item(14, 10).Value = 0.1884967751759062
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 3.057109
$ws.Cells.Item(14, 14).Value = 9.171327
$ws.Cells.Item(14, 15).Value = 0.02694952608666365
$ws.Cells.Item(14, 16).Value = 0.02694952608666365
$ws.Cells.Item(14, 17).Value = 3.591531395617
$ws.Cells.Item(14, 18).Value = 32.323782560553
$ws.Cells.Item(14, 19).Value = 0.005079898759855058
$ws.Cells.Item(14, 20).Value = 0.005079898759855057

# Row 15: sCs -> FAPs
$ws.Cells.Item(15, 1).Value = "sCs"
$ws.Cells.Item(15, 2).Value = "Lama3"
$ws.Cells.Item(15, 3).Value = "Sdc2"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 1.174813
$ws.Cells.Item(15, 8).Value = 3.524439
$ws.Cells.Item(15, 9).Value = 0.1884967751759062
$ws.Cells.Item(15, 10).Value = 0.1884967751759062
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 89.02756
$ws.Cells.Item(15, 14).Value = 267.08268
$ws.Cells.Item(15, 15).Value = 0.7848102735793893
$ws.Cells.Item(15, 16).Value = 0.7848102735793893
$ws.Cells.Item(15, 17).Value = 104.59073484628
$ws.Cells.Item(15, 18).Value = 941.3166136165199
$ws.Cells.Item(15, 19).Value = 0.1479342056946356
$ws.Cells.Item(15, 20).Value = 0.1479342056946356

# Row 16: sCs -> M2
$ws.Cells.Item(16, 1).Value = "sCs"
$ws.Cells.Item(16, 2).Value = "Lama3"
$ws.Cells.Item(16, 3).Value = "Sdc2"
$ws.Cells.Item(16, 4).Value = "M2"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 1.174813
$ws.Cells.Item(16, 8).Value = 3.524439
$ws.Cells.Item(16, 9).Value = 0.1884967751759062
$ws.Cells.Item(16, 10).Value = 0.1884967751759062
$ws.Cells.Item(16, 11).Value = 2
$ws.Cells.Item(16, 12).Value = 0.6666666666666666
$ws.Cells.Item(16, 13).Value = 0.184005
$ws.Cells.Item(16, 14).Value = 0.5520149999999999
$ws.Cells.Item(16, 15).Value = 0.0016220709001794
$ws.Cells.Item(16, 16).Value = 0.0016220709001794
$ws.Cells.Item(16, 17).Value = 0.216171466065
$ws.Cells.Item(16, 18).Value = 1.945543194585
$ws.Cells.Item(16, 19).Value = 0.0003057551337904961
$ws.Cells.Item(16, 20).Value = 0.000305755133790496

# Row 17: sCs -> sCs
$ws.Cells.Item(17, 1).Value = "sCs"
$ws.Cells.Item(17, 2).Value = "Lama3"
$ws.Cells.Item(17, 3).Value = "Sdc2"
$ws.Cells.Item(17, 4).Value = "sCs"
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 1.174813
$ws.Cells.Item(17, 8).Value = 3.524439
$ws.Cells.Item(17, 9).Value = 0.1884967751759062
$ws.Cells.Item(17, 10).Value = 0.1884967751759062
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 21.16964733333333
$ws.Cells.Item(17, 14).Value = 63.508942
$ws.Cells.Item(17, 15).Value = 0.1866181294337677
$ws.Cells.Item(17, 16).Value = 0.1866181294337677
$ws.Cells.Item(17, 17).Value = 24.87037689261534
$ws.Cells.Item(17, 18).Value = 223.833392033538
$ws.Cells.Item(17, 19).Value = 0.03517691558762508
$ws.Cells.Item(17, 20).Value = 0.03517691558762507

